$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose new Price text parses as a plain number (so Excel would store
# it as a numeric value rather than literal text, e.g. "333.20" -> 333.2,
# dropping the trailing zero, or introducing floating point noise like
# "0.574" -> 0.57399999999999995). For these rows we force the cell to Text
# format first so the literal digits survive exactly as in the source
# report. Prices that contain a thousands-separator dot (e.g.
# "51.505.22") or other non-numeric characters already round-trip as text
# with no extra help needed.
$TextPriceRows = @(5, 6, 8, 9, 10, 11, 12, 13, 17, 23, 24, 26, 30, 31, 32, 34, 37, 40, 41, 42, 43, 44, 51)

# Helper: set Price (D) and Volume(1h) (E) cells for a row
function Set-Row($row, $price, $volume) {
    if ($null -ne $price) {
        $cell = $ws.Cells.Item($row, 4)
        if ($TextPriceRows -contains $row) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $price
    }
    if ($null -ne $volume) {
        $ws.Cells.Item($row, 5).Value = "  $volume  "
    }
}

# Row 2 - Bitcoin
Set-Row 2 "51.505.22" "+3.75%"
# Row 3 - Ethereum
Set-Row 3 "2.759.46" "+4.50%"
# Row 4 - TetherUSD
Set-Row 4 $null "+0.03%"
# Row 5 - Solana
Set-Row 5 "116.25" "+2.45%"
# Row 6 - BNB
Set-Row 6 "333.20" "+2.65%"
# Row 7 - XRP
Set-Row 7 $null "+1.95%"
# Row 8 - USDC
Set-Row 8 "0.999" "-0.05%"
# Row 9 - Cardano
Set-Row 9 "0.574" "+5.16%"
# Row 10 - Avalanche
Set-Row 10 "41.74" "+4.34%"
# Row 11 - Dogecoin
Set-Row 11 "0.0868" "+6.62%"
# Row 12 - Chainlink
Set-Row 12 "20.30" "+2.38%"
# Row 13 - TRON
Set-Row 13 "0.130" "+2.29%"
# Row 14 - Polkadot
Set-Row 14 $null "+4.50%"
# Row 15 - WrappedliquidstakedEther2.0
Set-Row 15 "3.194.42" "+4.45%"
# Row 16 - WrappedEther
Set-Row 16 "2.764.24" "+4.72%"
# Row 17 - Polygon
Set-Row 17 "0.890" "+3.22%"
# Row 18 - WrappedBTC
Set-Row 18 "51.558.78" "+3.95%"
# Row 19 - ImmutableX
Set-Row 19 $null "+9.32%"
# Row 20 - InternetComputer(DFINITY)
Set-Row 20 $null "+4.18%"
# Row 21 - Uniswap
Set-Row 21 $null "+2.13%"
# Row 22 - ShibaInu
Set-Row 22 "0.0₃0977" "+2.98%"
# Row 23 - BitcoinCash
Set-Row 23 "278.55" "+3.10%"
# Row 24 - Litecoin
Set-Row 24 "69.66" "+0.96%"
# Row 25 - PancakeSwap
Set-Row 25 $null "+5.75%"
# Row 26 - EthereumClassic
Set-Row 26 "26.80" "+1.69%"
# Row 27 - Dai
Set-Row 27 $null "-0.05%"
# Row 28 - Cosmos
Set-Row 28 $null "-1.67%"
# Row 29 - Toncoin
Set-Row 29 $null "-0.18%"
# Row 30 - Kaspa
Set-Row 30 "0.142" "+2.09%"
# Row 31 - InjectiveProtocol
Set-Row 31 "35.08" "-0.53%"
# Row 32 - OKB
Set-Row 32 "50.06" "+0.97%"
# Row 33 - Filecoin
Set-Row 33 $null "+1.30%"
# Row 34 - Hedera
Set-Row 34 "0.0824" "+1.21%"
# Row 35 - FirstDigitalUSD
Set-Row 35 $null "-0.20%"
# Row 36 - Celestia
Set-Row 36 $null "-0.44%"
# Row 37 - RenderToken
Set-Row 37 "5.00" "+0.73%"
# Row 38 - ARBITRUM
Set-Row 38 $null "+1.44%"
# Row 39 - LidoDAOToken
Set-Row 39 $null "+3.45%"
# Row 40 - VeChain
Set-Row 40 "0.0353" "+7.94%"
# Row 41 - Monero
Set-Row 41 "127.40" "+0.30%"
# Row 42 - EnergySwap
Set-Row 42 "23.16" "+3.80%"

# Rows 43 & 44 - Stellar and WEMIXToken swap places (content swap, not just a move)
$ws.Cells.Item(43, 2).Value = "Stellar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-Row 43 "0.114" "+2.96%"

$ws.Cells.Item(44, 2).Value = "WEMIXToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-Row 44 "2.31" "+7.45%"

# Row 45 - Stacks
Set-Row 45 $null "+14.92%"
# Row 46 - Maker
Set-Row 46 "2.090.94" "+1.28%"
# Row 47 - NEARProtocol
Set-Row 47 $null "+2.36%"
# Row 48 - ApeXProtocol
Set-Row 48 $null "+4.00%"
# Row 49 - THORChain
Set-Row 49 $null "+5.66%"
# Row 50 - FraxShare
Set-Row 50 $null "+0.53%"
# Row 51 - MultiversX
Set-Row 51 "60.10" "+1.32%"

Write-Host "Applied cryptos update"
